# DPLKLib_Report.xlsx - 15/02/2023 / 16/02/2023 maintenance edits
# - Update the "URL" cell text on the Global sheet to the new server IP
# - Move the active selection to E2
# - Narrow column D to fit the shorter label

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# 1) Cell D2 holds the hyperlinked "URL" value - only the displayed text
#    changes (the hyperlink target itself is left as-is).
$ws.Range("D2").Value = "http://192.168.168.111/"

# 2) Narrow column D from its bestFit width down to 22 characters.
#    (21 + 5/6 round-trips through the host's column-width model to the
#    exact stored width of 22.)
$ws.Columns.Item(4).ColumnWidth = 21.1666666666667

# 3) Move the selection/active cell to E2.
$ws.Range("E2").Select()
